# Daily cryptos-list refresh (GitHub Actions bot): update Price/Volume(1h)
# columns from the latest coinranking.com scrape, and fix two coin rows
# whose rank order flipped (Avalanche/WrappedBTC, EnergySwap/Stacks).
#
# Numeric-looking Price strings (e.g. "7.45") must stay TEXT, matching the
# source inlineStr cells -- plain Range.Value assignment of such a string
# gets auto-coerced to a number, so for those cells we briefly force the
# cell to Text format, assign, then reset the style back to Normal so no
# stray number-format style is left behind on the cell.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "60.401.67"
$ws.Range("E2").Value = "  -5.81%  "
$ws.Range("D3").Value = "3.303.82"
$ws.Range("E3").Value = "  -5.03%  "
$ws.Range("E4").Value = "  +0.04%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "562.67"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -3.88%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "129.74"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -1.69%  "
$ws.Range("E7").Value = "  +0.02%  "
$ws.Range("D8").Value = "3.303.40"
$ws.Range("E8").Value = "  -5.01%  "
$ws.Range("E9").Value = "  -2.22%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "7.45"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -3.15%  "
$ws.Range("E11").Value = "  -4.94%  "
$ws.Range("E12").Value = "  -3.62%  "
$ws.Range("D13").Value = "3.857.54"
$ws.Range("E13").Value = "  -5.17%  "
$ws.Range("E14").Value = "  -0.31%  "
$ws.Range("D15").Value = "3.294.07"
$ws.Range("E15").Value = "  -5.30%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.0000166"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -5.99%  "
$ws.Range("B17").Value = "WrappedBTC"
$ws.Range("C17").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D17").Value = "60.581.02"
$ws.Range("E17").Value = "  -5.49%  "
$ws.Range("B18").Value = "Avalanche"
$ws.Range("C18").Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "24.19"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -4.06%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "5.68"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.07%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "13.29"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -0.85%  "
$ws.Range("E21").Value = "  -10.22%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "349.73"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -9.34%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.554"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -2.36%  "
$ws.Range("E24").Value = "  +0.16%  "
$ws.Range("D25").Value = "3.429.85"
$ws.Range("E25").Value = "  -5.16%  "
$ws.Range("E26").Value = "  -7.16%  "
$ws.Range("E27").Value = "  -3.60%  "
$ws.Range("E28").Value = "  -0.16%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "7.38"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +3.86%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.45"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +1.41%  "
$ws.Range("E31").Value = "  -2.03%  "
$ws.Range("E32").Value = "  -1.81%  "
$ws.Range("E33").Value = "  -5.63%  "
$ws.Range("D35").Value = "3.328.24"
$ws.Range("E35").Value = "  -5.02%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "22.60"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -1.64%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "5.34"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +2.25%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "6.76"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -0.12%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.48"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -1.52%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "157.13"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -2.86%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.0752"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -3.52%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.998"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -0.12%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "40.90"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -1.29%  "
$ws.Range("E44").Value = "  -0.17%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.742"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -7.01%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.17"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +2.85%  "
$ws.Range("B47").Value = "Stacks"
$ws.Range("C47").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.54"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -5.17%  "
$ws.Range("B48").Value = "EnergySwap"
$ws.Range("C48").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "22.54"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -3.87%  "
$ws.Range("E49").Value = "  -0.70%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "21.75"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +5.84%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.861"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -4.53%  "
